# [Experiment 2] Change description column in the queries.xlsx
#
# The "Description " header in F1 (shared with B1) is replaced with a
# single space, and the sheet view's scroll/selection is updated to
# reflect where the author was working (F10) when the change was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Replace the column F header text with a single space.
$ws.Range("F1").Value = " "

# Move the viewport/selection the way the author left it.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("F10").Select()
